$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF") ---
# Copy the formatting of the existing header cell H1 (bold, border, centered)
# onto I1 and J1 so they match the rest of the header row, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data columns I2:I60 and J2:J60 ---
$I = @(7, 9, 6, 6, 7, 6, 8, 8, 6, 5, 6, 5, 6, 6, 7, 7, 7, 8, 8, 5, 6, 8, 5, 7, 6, 6, 8, 9, 8, 7, 9, 7, 9, 8, 8, 9, 6, 8, 5, 6, 7, 9, 8, 7, 8, 6, 7, 6, 9, 5, 6, 9, 6, 6, 7, 5, 6, 5, 8)
$J = @(7, 9, 6, 7, 8, 7, 9, 9, 6, 5, 6, 5, 6, 6, 8, 8, 8, 8, 8, 5, 6, 8, 6, 7, 6, 6, 8, 9, 8, 7, 9, 7, 9, 9, 8, 9, 7, 8, 6, 6, 8, 9, 8, 7, 8, 6, 7, 6, 9, 5, 6, 9, 6, 6, 7, 6, 6, 5, 8)

for ($r = 2; $r -le 60; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $I[$idx]
    $ws.Cells.Item($r, 10).Value = $J[$idx]
}
